$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" summary text ---
$ws1 = $wb.Worksheets("Hoja1")

$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.58 = 50326.92 pesos`n✅ 50326.92 pesos = 12.52 = 980.0 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $conversionText

# --- Sheet "tasas": update the rate table values ---
$ws2 = $wb.Worksheets("tasas")

$ws2.Range("N10").Value = 79.5
$ws2.Range("O10").Value = 4000.99

$ws2.Range("N12").Value = 4021
$ws2.Range("O12").Value = 78.3
